# Append a new data row (row 48) to each of the four sensor-log sheets,
# mirroring the existing row layout (A: time, B: 总长, C: ID, D: 实际长度,
# E: 和校验, F: 总长_DEC, G: ID_DEC, H: 实际长度_DEC, I: 和校验_DEC).

$wb = $excel.ActiveWorkbook

$rowsToAdd = @(
    @{ Sheet = "ROW35-FE-LIFTER";  A = "2025-03-06 07:42:06"; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"; D = "0x01,0x90,"; E = "0x d"; F = 400; G = "568631262647113770877196"; H = 400; I = 13 },
    @{ Sheet = "ROW35-MID-LIFTER"; A = "2025-03-06 07:29:35"; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; D = "0x01,0x90,"; E = "0x e"; F = 400; G = "568631262647113770942732"; H = 400; I = 14 },
    @{ Sheet = "ROW02-FE-LIFTER";  A = "2025-03-06 07:51:45"; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"; D = "0x01,0x90,"; E = "0xff";  F = 400; G = "568631262647113769959692"; H = 400; I = 255 },
    @{ Sheet = "ROW02-MID-LIFTER"; A = "2025-03-06 07:41:15"; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x01,0x90,"; E = "0x 3"; F = 400; G = "568631262647113769959692"; H = 400; I = 3 }
)

foreach ($row in $rowsToAdd) {
    $ws = $wb.Worksheets.Item($row.Sheet)
    $r = 48

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F

    # Column G holds a 24-digit integer string. A plain assignment would be
    # auto-detected as a number and rounded to double precision, so force
    # text entry, then drop back to the default "Normal" style so the cell
    # doesn't end up with a stray explicit number format like its neighbours.
    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 7).Style = "Normal"

    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
}
